# gui create, select and run functions
# - Remove the obsolete "MENSAJE" column header (F1), which drops the
#   now-unused "MENSAJE" shared string from the workbook's string table.
# - Populate the sample RIF/code row (A2:B2) used by the new create/select/
#   run GUI functions.
# - Leave the selection on F22, matching where the user ended up in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "MENSAJE" header text from F1 (column stays, header text removed).
$ws.Range("F1").ClearContents()

# Fill in the sample RIF number and code values on row 2.
$ws.Range("A2").Value = 5524325124
$ws.Range("B2").Value = 233342

# Move/leave the active selection at F22.
$ws.Range("F22").Select()
